# Scheduled runner update: refresh market-price-derived leve profit figures
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
# LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns H:N) across all job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 500.4
$ws.Range("I18").Value = 251
$ws.Range("J18").Value = 666.6667
$ws.Range("K18").Value = 251
$ws.Range("L18").Value = 666.6667
$ws.Range("M18").Value = 33
$ws.Range("N18").Value = -1234.6667
$ws.Range("H76").Value = 4048.276
$ws.Range("J76").Value = 3107.6924
$ws.Range("L76").Value = 3107.6924
$ws.Range("N76").Value = -3737.6924
$ws.Range("H79").Value = 4048.276
$ws.Range("J79").Value = 3107.6924
$ws.Range("L79").Value = 3107.6924
$ws.Range("N79").Value = -5291.6924
$ws.Range("H80").Value = 2898.2666
$ws.Range("I80").Value = 955.2
$ws.Range("J80").Value = 6784.4
$ws.Range("K80").Value = 2865.6
$ws.Range("L80").Value = 20353.2
$ws.Range("M80").Value = -1867.6
$ws.Range("N80").Value = -22349.2
$ws.Range("H83").Value = 2898.2666
$ws.Range("I83").Value = 955.2
$ws.Range("J83").Value = 6784.4
$ws.Range("K83").Value = 8596.800000000001
$ws.Range("L83").Value = 61059.6
$ws.Range("M83").Value = -3604.800000000001
$ws.Range("N83").Value = -71043.60000000001
$ws.Range("H98").Value = 1341.375
$ws.Range("I98").Value = 954.1667
$ws.Range("J98").Value = 2503
$ws.Range("K98").Value = 954.1667
$ws.Range("L98").Value = 2503
$ws.Range("M98").Value = 543.8333
$ws.Range("N98").Value = -5499
$ws.Range("H113").Value = 2740
$ws.Range("I113").Value = 2700
$ws.Range("K113").Value = 2700
$ws.Range("M113").Value = 554
$ws.Range("H122").Value = 1341.375
$ws.Range("I122").Value = 954.1667
$ws.Range("J122").Value = 2503
$ws.Range("K122").Value = 2862.5001
$ws.Range("L122").Value = 7509
$ws.Range("M122").Value = -412.5001000000002
$ws.Range("N122").Value = -12409
$ws.Range("H133").Value = 60780
$ws.Range("J133").Value = 60780
$ws.Range("L133").Value = 60780
$ws.Range("N133").Value = -70900
$ws.Range("H135").Value = 1010.2727
$ws.Range("I135").Value = 1054.5
$ws.Range("K135").Value = 9490.5
$ws.Range("M135").Value = -6955.5
$ws.Range("H137").Value = 1276
$ws.Range("I137").Value = 539
$ws.Range("J137").Value = 2160.4
$ws.Range("K137").Value = 1617
$ws.Range("L137").Value = 6481.200000000001
$ws.Range("M137").Value = 933
$ws.Range("N137").Value = -11581.2

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 824280.0600000001
$ws.Range("I32").Value = 1114359
$ws.Range("J32").Value = 18505.334
$ws.Range("K32").Value = 1114359
$ws.Range("L32").Value = 18505.334
$ws.Range("M32").Value = -1114072
$ws.Range("N32").Value = -19079.334
$ws.Range("H61").Value = 6668554.5
$ws.Range("I61").Value = 10418327
$ws.Range("J61").Value = 2291.2778
$ws.Range("K61").Value = 10418327
$ws.Range("L61").Value = 2291.2778
$ws.Range("M61").Value = -10418115
$ws.Range("N61").Value = -2715.2778
$ws.Range("H63").Value = 3544.3333
$ws.Range("I63").Value = 2885.5
$ws.Range("J63").Value = 4143.273
$ws.Range("K63").Value = 2885.5
$ws.Range("L63").Value = 4143.273
$ws.Range("M63").Value = -2199.5
$ws.Range("N63").Value = -5515.273
$ws.Range("H66").Value = 3544.3333
$ws.Range("I66").Value = 2885.5
$ws.Range("J66").Value = 4143.273
$ws.Range("K66").Value = 14427.5
$ws.Range("L66").Value = 20716.365
$ws.Range("M66").Value = -10995.5
$ws.Range("N66").Value = -27580.365
$ws.Range("H132").Value = 4827.436
$ws.Range("I132").Value = 4981.5835
$ws.Range("J132").Value = 4580.8
$ws.Range("K132").Value = 14944.7505
$ws.Range("L132").Value = 13742.4
$ws.Range("M132").Value = -12414.7505
$ws.Range("N132").Value = -18802.4
$ws.Range("H136").Value = 6668554.5
$ws.Range("I136").Value = 10418327
$ws.Range("J136").Value = 2291.2778
$ws.Range("K136").Value = 31254981
$ws.Range("L136").Value = 6873.8334
$ws.Range("M136").Value = -31252431
$ws.Range("N136").Value = -11973.8334

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 11505.333
$ws.Range("I24").Value = 1016
$ws.Range("J24").Value = 16750
$ws.Range("K24").Value = 1016
$ws.Range("L24").Value = 16750
$ws.Range("M24").Value = -781
$ws.Range("N24").Value = -17220
$ws.Range("H86").Value = 2513.5
$ws.Range("I86").Value = 2513.5
$ws.Range("K86").Value = 2513.5
$ws.Range("M86").Value = -1390.5
$ws.Range("H89").Value = 2513.5
$ws.Range("I89").Value = 2513.5
$ws.Range("K89").Value = 12567.5
$ws.Range("M89").Value = -6951.5
$ws.Range("H105").Value = 2450.1667
$ws.Range("I105").Value = 2140.2
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 2140.2
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -393.1999999999998
$ws.Range("N105").Value = -7494
$ws.Range("H134").Value = 2786.5527
$ws.Range("I134").Value = 2883.04
$ws.Range("J134").Value = 2601
$ws.Range("K134").Value = 8649.119999999999
$ws.Range("L134").Value = 7803
$ws.Range("M134").Value = -6114.119999999999
$ws.Range("N134").Value = -12873

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3371.1538
$ws.Range("I62").Value = 3402.0833
$ws.Range("K62").Value = 3402.0833
$ws.Range("M62").Value = -2778.0833
$ws.Range("H65").Value = 3371.1538
$ws.Range("I65").Value = 3402.0833
$ws.Range("K65").Value = 17010.4165
$ws.Range("M65").Value = -13890.4165
$ws.Range("H132").Value = 8773951
$ws.Range("I132").Value = 1715.4546
$ws.Range("J132").Value = 20835774
$ws.Range("K132").Value = 5146.3638
$ws.Range("L132").Value = 62507322
$ws.Range("M132").Value = -2616.3638
$ws.Range("N132").Value = -62512382

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 932.5208
$ws.Range("I68").Value = 659
$ws.Range("J68").Value = 1004.5
$ws.Range("K68").Value = 1977
$ws.Range("L68").Value = 3013.5
$ws.Range("M68").Value = -1166
$ws.Range("N68").Value = -4635.5
$ws.Range("H71").Value = 932.5208
$ws.Range("I71").Value = 659
$ws.Range("J71").Value = 1004.5
$ws.Range("K71").Value = 5931
$ws.Range("L71").Value = 9040.5
$ws.Range("M71").Value = -1875
$ws.Range("N71").Value = -17152.5
$ws.Range("H76").Value = 2979.0425
$ws.Range("J76").Value = 3000.3262
$ws.Range("L76").Value = 9000.9786
$ws.Range("N76").Value = -9766.9786
$ws.Range("H79").Value = 2979.0425
$ws.Range("J79").Value = 3000.3262
$ws.Range("L79").Value = 9000.9786
$ws.Range("N79").Value = -11652.9786

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H122").Value = 1346.8
$ws.Range("I122").Value = 1263.2
$ws.Range("J122").Value = 1514
$ws.Range("K122").Value = 3789.6
$ws.Range("L122").Value = 4542
$ws.Range("M122").Value = -1339.6
$ws.Range("N122").Value = -9442

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12800.777
$ws.Range("I22").Value = 933.3333
$ws.Range("K22").Value = 933.3333
$ws.Range("M22").Value = -638.3333
$ws.Range("H27").Value = 12800.777
$ws.Range("I27").Value = 933.3333
$ws.Range("K27").Value = 933.3333
$ws.Range("M27").Value = -826.3333
$ws.Range("H132").Value = 3527.1707
$ws.Range("I132").Value = 3162.1538
$ws.Range("J132").Value = 4159.8667
$ws.Range("K132").Value = 9486.4614
$ws.Range("L132").Value = 12479.6001
$ws.Range("M132").Value = -6956.4614
$ws.Range("N132").Value = -17539.6001

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 397.75
$ws.Range("I100").Value = 277.33334
$ws.Range("J100").Value = 470
$ws.Range("K100").Value = 554.66668
$ws.Range("L100").Value = 940
$ws.Range("M100").Value = -13.66668000000004
$ws.Range("N100").Value = -2022

Write-Output "Applied scheduled market-data refresh across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
